$wb = $excel.ActiveWorkbook

# --- Sheet "dimenzija za m=4": B1 currently holds text "4", convert to real number 4 ---
$ws1 = $wb.Worksheets.Item("dimenzija za m=4")
$ws1.Range("B1").Value = 4

# --- Sheet "dimenzija za m=5-10": clear the "Ni definirano za m >= n" placeholder text ---
# Cells keep existing (as empty cells) rather than being removed outright, so a
# no-op border touch follows the value-clear to keep them present in sheetData
# (matches the original file, which kept them as empty inlineStr cells).
$ws2 = $wb.Worksheets.Item("dimenzija za m=5-10")
$cellsToClear = "D2", "E2", "F2", "G2", "E3", "F3", "G3", "F4", "G4", "G5"
foreach ($ref in $cellsToClear) {
    $cell = $ws2.Range($ref)
    $cell.Value = ""
    $cell.Borders.LineStyle = -4142
}
